# Add a new "Project Status" column (R) to Sheet1, classifying each
# project row as either "In Progress" or "Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("R1").Value = "Project Status"

# Row ranges (inclusive) that should be marked "In Progress"
$inProgressRanges = @(
    @(2, 31),
    @(46, 50),
    @(64, 68),
    @(75, 86)
)

# Row ranges (inclusive) that should be marked "Done"
$doneRanges = @(
    @(32, 45),
    @(51, 63),
    @(69, 74)
)

foreach ($range in $inProgressRanges) {
    $startRow = $range[0]
    $endRow = $range[1]
    $ws.Range("R$startRow`:R$endRow").Value = "In Progress"
}

foreach ($range in $doneRanges) {
    $startRow = $range[0]
    $endRow = $range[1]
    $ws.Range("R$startRow`:R$endRow").Value = "Done"
}

# Update the visible selection to match the saved view state
$ws.Range("U74").Select()
